$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.6402153558052435
$ws.Range("D2").Value = 0.6305742837354592
$ws.Range("E2").Value = 0.2764798242788156
$ws.Range("F2").Value = 0.5157564879398394
$ws.Range("G2").Value = 0.5004541326067211
$ws.Range("H2").Value = 0.624727884425094
$ws.Range("I2").Value = 0.2453740332067915
$ws.Range("J2").Value = 0.4568520167462022

$ws.Range("C3").Value = 0.7652153558052435
$ws.Range("D3").Value = 0.7285968706271466
$ws.Range("E3").Value = 0.3782158706398598
$ws.Range("F3").Value = 0.6240093656907499
$ws.Range("G3").Value = 0.485921889191644
$ws.Range("H3").Value = 0.722343162477736
$ws.Range("I3").Value = 0.3270598505952984
$ws.Range("J3").Value = 0.5117749674215595

$ws.Range("C4").Value = 0.8886938202247191
$ws.Range("D4").Value = 0.7496289700207777
$ws.Range("E4").Value = 0.4736095964292549
$ws.Range("F4").Value = 0.7039774622249172
$ws.Range("G4").Value = 0.5140781108083561
$ws.Range("H4").Value = 0.7346625766871165
$ws.Range("I4").Value = 0.3972993705211142
$ws.Range("J4").Value = 0.5486800193388622

$ws.Range("C5").Value = 0.9046114232209738
$ws.Range("D5").Value = 0.7607457349220483
$ws.Range("E5").Value = 0.4996195511294804
$ws.Range("F5").Value = 0.7216589030908341
$ws.Range("G5").Value = 0.5059037238873751
$ws.Range("H5").Value = 0.7424797150207797
$ws.Range("I5").Value = 0.3815912287648544
$ws.Range("J5").Value = 0.5433248892243364

$ws.Range("C6").Value = 0.9702715355805244
$ws.Range("D6").Value = 0.7621097117980467
$ws.Range("E6").Value = 0.5650893777043705
$ws.Range("F6").Value = 0.7658235416943139
$ws.Range("G6").Value = 0.5059037238873751
$ws.Range("H6").Value = 0.7403522659806056
$ws.Range("I6").Value = 0.3962485499877949
$ws.Range("J6").Value = 0.5475015132852584

$ws.Range("C7").Value = 0.9731975655430711
$ws.Range("D7").Value = 0.7770710540078305
$ws.Range("E7").Value = 0.6177959043413878
$ws.Range("F7").Value = 0.7893548412974298
$ws.Range("G7").Value = 0.5095367847411444
$ws.Range("H7").Value = 0.7522758757173956
$ws.Range("I7").Value = 0.4126899351101684
$ws.Range("J7").Value = 0.5581675318562361

$ws.Range("C8").Value = 0.9860720973782772
$ws.Range("D8").Value = 0.7789368047604913
$ws.Range("E8").Value = 0.6670591584846525
$ws.Range("F8").Value = 0.8106893535411404
$ws.Range("G8").Value = 0.5122615803814714
$ws.Range("H8").Value = 0.7485652087868593
$ws.Range("I8").Value = 0.4198586675077124
$ws.Range("J8").Value = 0.5602284855586811

$ws.Range("C9").Value = 0.9832631086142322
$ws.Range("D9").Value = 0.7890288201953385
$ws.Range("E9").Value = 0.7072696882658795
$ws.Range("F9").Value = 0.8265205390251501
$ws.Range("G9").Value = 0.5095367847411444
$ws.Range("H9").Value = 0.7558875915297842
$ws.Range("I9").Value = 0.4328464779303688
$ws.Range("J9").Value = 0.5660902847337658

$ws.Range("C10").Value = 0.9928604868913857
$ws.Range("D10").Value = 0.7932055576757269
$ws.Range("E10").Value = 0.7556263227652499
$ws.Range("F10").Value = 0.8472307891107875
$ws.Range("G10").Value = 0.5204359673024523
$ws.Range("H10").Value = 0.7571739560657036
$ws.Range("I10").Value = 0.4568105041770552
$ws.Range("J10").Value = 0.5781401425150704

$ws.Range("C11").Value = 0.9952013108614233
$ws.Range("D11").Value = 0.8005201486946812
$ws.Range("E11").Value = 0.7836565374349646
$ws.Range("F11").Value = 0.8597926656636897
$ws.Range("G11").Value = 0.5277020890099909
$ws.Range("H11").Value = 0.759499307342173
$ws.Range("I11").Value = 0.4565078155115433
$ws.Range("J11").Value = 0.5812364039545691
